$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.978.07"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "'2.206.08"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'230.52"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("D7").Value = "'60.55"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").Value = "'0.0897"
$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "'2.533.50"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("D13").Value = "'15.37"
$ws.Range("E13").Value = "  -1.64%  "

$ws.Range("D14").Value = "'21.92"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").Value = "'0.794"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "'5.56"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'2.202.06"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").Value = "'41.891.89"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").Value = "'0.0₃0932"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("D21").Value = "'71.77"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").Value = "'242.47"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("E25").Value = "  +2.37%  "

$ws.Range("D26").Value = "'9.59"
$ws.Range("E26").Value = "  +0.99%  "

$ws.Range("D27").Value = "'168.73"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").Value = "'20.24"
$ws.Range("E29").Value = "  +1.81%  "

$ws.Range("D30").Value = "'1.43"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").Value = "'4.92"
$ws.Range("E33").Value = "  -2.45%  "

$ws.Range("D34").Value = "'4.58"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("D35").Value = "'0.0644"
$ws.Range("E35").Value = "  +3.69%  "

$ws.Range("D36").Value = "'6.28"
$ws.Range("E36").Value = "  -4.83%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.32"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  -4.55%  "

$ws.Range("D39").Value = "'0.0248"
$ws.Range("E39").Value = "  +6.16%  "

$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("D41").Value = "'8.54"
$ws.Range("E41").Value = "  -1.58%  "

$ws.Range("D42").Value = "'0.000224"
$ws.Range("E42").Value = "  -5.53%  "

$ws.Range("D43").Value = "'0.0951"
$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("E44").Value = "  +1.34%  "

$ws.Range("D45").Value = "'96.50"
$ws.Range("E45").Value = "  -2.58%  "

$ws.Range("D46").Value = "'1.455.87"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").Value = "'4.28"
$ws.Range("E47").Value = "  -12.33%  "

$ws.Range("D48").Value = "'2.73"
$ws.Range("E48").Value = "  -1.76%  "

$ws.Range("D49").Value = "'16.03"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  -1.58%  "

$ws.Range("E51").Value = "  +2.02%  "
